$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 174
$ws.Range("A2").Value = 50
$ws.Range("A3").Value = 116
$ws.Range("A4").Value = 149
$ws.Range("A5").Value = 34
$ws.Range("A6").Value = 99
$ws.Range("A7").Value = 81
$ws.Range("A8").Value = 150
